$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 57.71561539619783
$ws.Range("B3").Value = 10.75824539048188
$ws.Range("B4").Value = 10.05735555675508
$ws.Range("B5").Value = 8.701212397762976
$ws.Range("B6").Value = 6.738617788281282
$ws.Range("B7").Value = 6.02895347052096
